$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test-case table, entered roughly in the order the author actually typed
# it (reconstructed from the shared-string insertion order): the row 2
# header was typed B2, A2, D2, then they jumped down to A3 before coming
# back to finish C2 - everything after that is a plain left-to-right,
# top-to-bottom fill.
$ws.Range("B2").Value = "Color Plane "
$ws.Range("A2").Value = "Input: Name:"
$ws.Range("D2").Value = "Output: Destination"

$ws.Range("A3").Value = "Theresa "

$ws.Range("C2").Value = "Output: Plane Class"

$ws.Range("B3").Value = "pink"
$ws.Range("C3").Value = "first"
$ws.Range("D3").Value = "costa rica "

$ws.Range("A4").Value = "Jill"
$ws.Range("B4").Value = "blue"
$ws.Range("C4").Value = "economy"
$ws.Range("D4").Value = "south africa "

$ws.Range("A5").Value = "Kiley"
$ws.Range("B5").Value = "purple"
$ws.Range("C5").Value = "business"
$ws.Range("D5").Value = "spain "

$ws.Range("A6").Value = "Theresa "
$ws.Range("B6").Value = "purple"
$ws.Range("C6").Value = "first "
$ws.Range("D6").Value = "spain "

$ws.Range("A7").Value = "Julianna "
$ws.Range("B7").Value = "blue"
$ws.Range("C7").Value = "first"
$ws.Range("D7").Value = "south africa "

# Column widths (ColumnWidth in chars; Excel stores raw XML width as
# ColumnWidth + 5/MaximumDigitWidth, so we back that padding out here to
# land on the target stored widths of 16.5 / 18 / 25.1640625)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(3).ColumnWidth = 17.166666666666668
$ws.Columns.Item(4).ColumnWidth = 24.330729166666668

# Selection state as left by the author
$ws.Range("E17:F17").Select()
